# "Generate Report for Handoff"
#
# The localization status report is regenerated: the two tracked files
# (c60e3de8-...md and c8c535ca-...md) swap logical "slots" in each sheet's
# 2-row table (row 2 / row 3), and the file that is now in row 3
# (c60e3de8-...md) moves from "Handed back: in sync with en-US" to
# "Ready for handoff" with fresh handoff timestamps and (for zh-cn/de-de)
# a stale-handback error detail message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$ov.Range("B2").Value = "e2e\c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"

$ov.Range("A3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$ov.Range("B3").Value = "e2e\c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-02 02:56:30"

foreach ($h in $ov.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$B`$2") {
        $h.TextToDisplay = "e2e\c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
    } elseif ($addr -eq "`$B`$3") {
        $h.TextToDisplay = "e2e\c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$zh.Range("G2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.zh-cn.xlf"
$zh.Range("I2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$zh.Range("J2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.zh-cn.xlf"

$zh.Range("A3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-02 02:56:26"
$zh.Range("I3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$zh.Range("J3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83aed6d0d3c6207ad315003c6446dd872a253475/e2e/c60e3de8-f0b1-463f-83d6-957c38bb26a9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7566e23389698541e2e92a6946ad9a4293ea7d5e/e2e/c60e3de8-f0b1-463f-83d6-957c38bb26a9.md."

foreach ($h in $zh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2" -or $addr -eq "`$I`$2") {
        $h.TextToDisplay = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
    } elseif ($addr -eq "`$A`$3" -or $addr -eq "`$I`$3") {
        $h.TextToDisplay = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
    }
}

$zh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$de.Range("G2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.de-de.xlf"
$de.Range("I2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
$de.Range("J2").Value = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.206cdaf7b6d6b8c70b547477a8ef777d7d03ebc5.de-de.xlf"

$de.Range("A3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.de-de.xlf"
$de.Range("H3").Value = "2016-09-02 02:56:30"
$de.Range("I3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
$de.Range("J3").Value = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.6f9abbd882ba195e39018fd4a0cc57668314096a.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83aed6d0d3c6207ad315003c6446dd872a253475/e2e/c60e3de8-f0b1-463f-83d6-957c38bb26a9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7566e23389698541e2e92a6946ad9a4293ea7d5e/e2e/c60e3de8-f0b1-463f-83d6-957c38bb26a9.md."

foreach ($h in $de.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2" -or $addr -eq "`$I`$2") {
        $h.TextToDisplay = "c8c535ca-e5aa-49e9-b056-1410c2db1b01.md"
    } elseif ($addr -eq "`$A`$3" -or $addr -eq "`$I`$3") {
        $h.TextToDisplay = "c60e3de8-f0b1-463f-83d6-957c38bb26a9.md"
    }
}

$de.Columns.Item(16).ColumnWidth = 39.17
